# Added NBER IMR draft to website
# - Replace the "Unconditional Cash Transfers and Child Mortality" working paper
#   entry (row 4 on "Working Papers") with the updated "Can Cash Transfers Save
#   Lives?" NBER draft, including a new abstract and a link to the new PDF.
# - Mark the Free-Riding / New Product Adoption entry on "Work in Progress" as
#   "Data collection in process" (italic lead-in) and add its project abstract.

$wb = $excel.ActiveWorkbook

$wsWorking = $wb.Worksheets.Item(1)      # "Working Papers"
$wsWip     = $wb.Worksheets.Item(3)      # "Work in Progress"

# ---------------------------------------------------------------------------
# 1. Working Papers sheet, row 4: new cash-transfer / mortality paper
# ---------------------------------------------------------------------------

$newTitle = @'
Can Cash Transfers Save Lives? Evidence from a Large-Scale Experiment in Kenya
'@

$newAbstract = @'
We estimate the impacts of large-scale unconditional cash transfers on child survival. One-time transfers of USD 1000 were provided to over 10,500 poor households across 653 randomized villages in Kenya. We collected census data on over 100,000 births, including on mortality and cause of death, and detailed data on household health behaviors. Unconditional cash transfers (accounting for spillovers) lead to 48\% fewer infant deaths before age one and 45\% fewer child deaths before age five. Detailed data on cause of death, transfer timing relative to birth, and the location of health facilities indicate that unconditional cash transfers and access to delivery care are complements in generating mortality reductions: the largest gains are estimated in neonatal and maternal causes of death largely preventable by appropriate obstetric care and among households living close to physician-staffed facilities and those who receive the transfer around the time of birth, and treatment leads to a large increase in hospital deliveries (by 45\%). The infant and child mortality declines are concentrated among poorer households with below median assets or predicted consumption. The transfers also result in a substantial decline of 51\% in female labor supply in the three months before and the three months after a birth, and improved child nutrition. Infant and child mortality largely revert to pre-program levels after cash transfers end. Despite not being the main aim of the original program, we show that unconditional cash transfers in this setting may be a cost-effective way to reduce infant and child deaths.
'@

$newLink = @'
https://github.com/gkilleen33/gkilleen33.github.io/blob/master/papers/working/GE-IMR.pdf
'@

# Title, authors (unchanged) and abstract text
$wsWorking.Range("A4").Value = $newTitle
$wsWorking.Range("C4").Value = $newAbstract

# Link cell + hyperlink (match the look of the other paper links in D2/D3)
$wsWorking.Range("D4").Value = $newLink
$wsWorking.Hyperlinks.Add($wsWorking.Range("D4"), $newLink)
$wsWorking.Range("D2").Copy()
$wsWorking.Range("D4").PasteSpecial(-4122)  # xlPasteFormats, reuse D2/D3 style
$excel.CutCopyMode = $false

# Latex flag column, same as the other fully-populated rows
$wsWorking.Range("E2").Copy()
$wsWorking.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsWorking.Range("E4").Value = 1

# ---------------------------------------------------------------------------
# 2. Work in Progress sheet, row 2: Free-Riding / New Product Adoption paper
# ---------------------------------------------------------------------------

$dataCollectionRun = @'
Data collection in process 
'@

$authorsRun = @'
with Luisa Cefala, Rédempteur Ntawiratsa and Nicholas Swanson
'@

$wipAbstract = @'
In low and middle-income (LMIC) countries, businesses often use technologies inside the frontier, innovate less, and slowly adopt new products and technologies (Cirera et al. 2022). The reasons for these facts are not well understood, particularly for small and microenterprises (e.g. Atkin et al. 2017). 
We investigate whether a lack of institutions to protect the value of intellectual property contributes to these facts. High-income economies tend to have strong patent systems to promote discoveries, and regulators permit exclusive dealing in retail environments due to similar forces. But LMIC firms are often informal and undifferentiated, meaning neighbours are likely to adopt discoveries of their competitors without compensation. 
This project focuses on this problem in the case of retail firms’ decision to adopt a new product-- a setting where firms face risk ex-ante because they do not know if demand will exist, but, ex-post, competitors can learn if demand is high by observing the first mover. We examine whether offering retailers exclusive access to supply of a new product promotes adoption. Our study also tests whether firms are colluding so that if there is a null result we can understand if it is due to collusion, a possible upside to the finding that some markets in LMICs are uncompetitive (Bergquist and Dinerstein 2020). 
'@

$combinedAuthorsCell = $dataCollectionRun + $authorsRun
$runSplit = $dataCollectionRun.Length

$wsWip.Range("B2").Value = $combinedAuthorsCell
$charsItalic = $wsWip.Range("B2").Characters(1, $runSplit)
$charsItalic.Font.Italic = $true
$charsNormal = $wsWip.Range("B2").Characters($runSplit + 1, $authorsRun.Length)
$charsNormal.Font.Italic = $false

$wsWip.Range("C2").Value = $wipAbstract

# Row grows to hold the full abstract, same cap used elsewhere in the workbook
$wsWip.Rows.Item(2).RowHeight = 409.5

# ---------------------------------------------------------------------------
# 3. Selection / view bookkeeping to mirror where the author left the cursor
# ---------------------------------------------------------------------------

$wsWip.Range("D2").Select()

$wsWorking.Activate()
$wsWorking.Range("B7").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
